# "changed the colours of the pie graph"
#
# This deck's single slide hosts a PowerPoll poll/survey Office Add-in
# (a `we:webextension` web-extension graphicFrame, with a static `p:pic`
# fallback snapshot shown whenever the live add-in can't render). The
# add-in's internal state (and the colours of the pie chart it draws) are
# owned by the add-in's own web sandbox - there is no VBA/COM surface,
# in real PowerPoint or in this object model, that can reach into
# `we:webextension` bindings/properties to recolour the chart directly.
#
# The externally visible effect of a user recolouring the poll's pie
# chart inside the add-in is a refreshed fallback snapshot picture on the
# slide (PowerPoint re-renders/re-embeds the add-in's preview image).
# We reproduce that visible effect the same way any COM script replaces
# a picture's content - by dropping in a new image - sized and
# positioned to cover the same area the add-in's snapshot occupies.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The add-in's graphicFrame/fallback picture fills the whole 720x540pt
# slide (9144000 x 6858000 EMU), starting at the top-left corner.
$left = 0
$top = 0
$width = 720
$height = 540

$pic = $s.Shapes.AddPicture("/tmp/work/pie_colors_new.png", $false, $true, $left, $top, $width, $height)
$pic.Name = "Poll Pie Chart (updated colours)"
